# Customer, Vendor Sample file chgs
# Applies the data updates described by the commit to the sample
# customer.xlsx workbook (sheets: Customers, State, Address, Sales, Currency).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Customers"
# ---------------------------------------------------------------------------
$wsCustomers = $wb.Worksheets.Item("Customers")
$wsCustomers.Range("A2").Value = "ABCD PVT LTD"
$wsCustomers.Range("B2").Value = "29GGGGG1314R9Z6"
$wsCustomers.Range("C2").Value = "ABCDE1234F"

# ---------------------------------------------------------------------------
# Sheet "State"
# ---------------------------------------------------------------------------
$wsState = $wb.Worksheets.Item("State")
$wsState.Range("E2").Value = "JOE"
$wsState.Range("F2").Value = 9876543210
$wsState.Range("G2").Value = "joe@gmail.com"
$wsState.Range("H2").Value = "ABCD PVT LTD"

$wsState.Range("E3").Value = "ALICE"
$wsState.Range("F3").Value = 9876543210
$wsState.Range("G3").Value = "alice@gmail.com"
$wsState.Range("H3").Value = "ABCD PVT LTD"

# ---------------------------------------------------------------------------
# Sheet "Address"
# ---------------------------------------------------------------------------
$wsAddress = $wb.Worksheets.Item("Address")
$wsAddress.Range("C2").Value = "29GGGGG1314R9Z6"
$wsAddress.Range("D2").Value = "CHENNAI"
$wsAddress.Range("F2").Value = "Abc"
$wsAddress.Range("G2").Value = "Defg"
$wsAddress.Range("H2").Value = "Hij"
$wsAddress.Range("I2").Value = 100001
$wsAddress.Range("J2").Value = 9876543210
$wsAddress.Range("K2").Value = "ABCD PVT LTD"

$wsAddress.Range("C3").Value = "29GGGGG1314R9Z6"
$wsAddress.Range("D3").Value = "CHENNAI"
$wsAddress.Range("F3").Value = "Abc"
$wsAddress.Range("G3").Value = "Defg"
$wsAddress.Range("H3").Value = "Hij"
$wsAddress.Range("I3").Value = 100001
$wsAddress.Range("J3").Value = 9876543210
$wsAddress.Range("K3").Value = "ABCD PVT LTD"

# ---------------------------------------------------------------------------
# Sheet "Sales" - header used to sit on row 2 (row 1 was blank) with data on
# rows 3-4; now the blank row is removed so the header moves to row 1 and the
# two data rows shift up to rows 2-3.
# ---------------------------------------------------------------------------
$wsSales = $wb.Worksheets.Item("Sales")
$wsSales.Rows.Item(1).Delete()

$wsSales.Range("A2").Value = "ABC"
$wsSales.Range("F2").Value = "ABCD PVT LTD"

$wsSales.Range("A3").Value = "XYZ"
$wsSales.Range("F3").Value = "ABCD PVT LTD"

# ---------------------------------------------------------------------------
# Sheet "Currency"
# ---------------------------------------------------------------------------
$wsCurrency = $wb.Worksheets.Item("Currency")
$wsCurrency.Range("B2").Value = "ABCD PVT LTD"
$wsCurrency.Rows.Item(3).Delete()

# ---------------------------------------------------------------------------
# Selections / active sheet, mirroring the saved view state in the workbook.
# ---------------------------------------------------------------------------
$wsCustomers.Range("F14").Select()
$wsState.Range("F10").Select()
$wsAddress.Range("K8").Select()
$wsSales.Range("F16").Select()
$wsCurrency.Range("D9").Select()

$wsState.Activate()
